$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 4302.923
$ws.Range("I113").Value = 5166.2856
$ws.Range("J113").Value = 3295.6667
$ws.Range("K113").Value = 5166.2856
$ws.Range("L113").Value = 3295.6667
$ws.Range("M113").Value = -1912.2856
$ws.Range("N113").Value = -9803.6667
$ws.Range("H132").Value = 1997
$ws.Range("I132").Value = 1697.4884
$ws.Range("K132").Value = 5092.4652
$ws.Range("M132").Value = -2562.4652

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1672.258
$ws.Range("J2").Value = 987
$ws.Range("L2").Value = 987
$ws.Range("N2").Value = -1213
$ws.Range("H37").Value = 42999
$ws.Range("I37").Value = 39999
$ws.Range("J37").Value = 44499
$ws.Range("K37").Value = 39999
$ws.Range("L37").Value = 44499
$ws.Range("M37").Value = -39726
$ws.Range("N37").Value = -45045
$ws.Range("H45").Value = 1486.9231
$ws.Range("I45").Value = 1194.25
$ws.Range("J45").Value = 4999
$ws.Range("K45").Value = 1194.25
$ws.Range("L45").Value = 4999
$ws.Range("M45").Value = -817.25
$ws.Range("N45").Value = -5753
$ws.Range("H55").Value = 39049
$ws.Range("J55").Value = 0
$ws.Range("L55").Value = 0
$ws.Range("N55").ClearContents()
$ws.Range("H74").Value = 54197.94
$ws.Range("I74").Value = 30159.064
$ws.Range("K74").Value = 30159.064
$ws.Range("M74").Value = -29285.064
$ws.Range("H77").Value = 54197.94
$ws.Range("I77").Value = 30159.064
$ws.Range("K77").Value = 150795.32
$ws.Range("M77").Value = -146427.32
$ws.Range("H80").Value = 39991.5
$ws.Range("J80").Value = 39991.5
$ws.Range("L80").Value = 39991.5
$ws.Range("N80").Value = -41987.5
$ws.Range("H83").Value = 39991.5
$ws.Range("J83").Value = 39991.5
$ws.Range("L83").Value = 119974.5
$ws.Range("N83").Value = -129958.5
$ws.Range("H116").Value = 1672.258
$ws.Range("J116").Value = 987
$ws.Range("L116").Value = 987
$ws.Range("N116").Value = -5575
$ws.Range("H132").Value = 6525.035
$ws.Range("I132").Value = 6932
$ws.Range("K132").Value = 20796
$ws.Range("M132").Value = -18266

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1672.258
$ws.Range("J3").Value = 987
$ws.Range("L3").Value = 987
$ws.Range("N3").Value = -1215
$ws.Range("H82").Value = 12254.1875
$ws.Range("I82").Value = 3015.5
$ws.Range("J82").Value = 39970.25
$ws.Range("K82").Value = 3015.5
$ws.Range("L82").Value = 39970.25
$ws.Range("M82").Value = -2632.5
$ws.Range("N82").Value = -40736.25
$ws.Range("H85").Value = 12254.1875
$ws.Range("I85").Value = 3015.5
$ws.Range("J85").Value = 39970.25
$ws.Range("K85").Value = 3015.5
$ws.Range("L85").Value = 39970.25
$ws.Range("M85").Value = -1689.5
$ws.Range("N85").Value = -42622.25
$ws.Range("H134").Value = 3576.9443
$ws.Range("I134").Value = 2892.3333
$ws.Range("K134").Value = 8676.999899999999
$ws.Range("M134").Value = -6141.999899999999

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 945.6667
$ws.Range("I16").Value = 349.5
$ws.Range("K16").Value = 349.5
$ws.Range("M16").Value = -62.5
$ws.Range("H31").Value = 2624.9614
$ws.Range("I31").Value = 1694
$ws.Range("K31").Value = 1694
$ws.Range("M31").Value = -1399
$ws.Range("H34").Value = 2624.9614
$ws.Range("I34").Value = 1694
$ws.Range("K34").Value = 1694
$ws.Range("M34").Value = -1492
$ws.Range("H50").Value = 69418.125
$ws.Range("I50").Value = 0
$ws.Range("K50").Value = 0
$ws.Range("M50").ClearContents()
$ws.Range("H59").Value = 19916.334
$ws.Range("I59").Value = 2375
$ws.Range("J59").Value = 54999
$ws.Range("K59").Value = 2375
$ws.Range("L59").Value = 54999
$ws.Range("M59").Value = -1230
$ws.Range("N59").Value = -57289
$ws.Range("H60").Value = 23947.637
$ws.Range("I60").Value = 9061.143
$ws.Range("J60").Value = 49999
$ws.Range("K60").Value = 9061.143
$ws.Range("L60").Value = 49999
$ws.Range("M60").Value = -8550.143
$ws.Range("N60").Value = -51021
$ws.Range("H68").Value = 62304.2
$ws.Range("J68").Value = 67880.25
$ws.Range("L68").Value = 67880.25
$ws.Range("N68").Value = -69378.25
$ws.Range("H71").Value = 62304.2
$ws.Range("J71").Value = 67880.25
$ws.Range("L71").Value = 203640.75
$ws.Range("N71").Value = -211128.75
$ws.Range("H99").Value = 4501.2
$ws.Range("I99").Value = 4754.4443
$ws.Range("J99").Value = 4121.3335
$ws.Range("K99").Value = 4754.4443
$ws.Range("L99").Value = 4121.3335
$ws.Range("M99").Value = -3256.4443
$ws.Range("N99").Value = -7117.3335
$ws.Range("H113").Value = 945.6667
$ws.Range("I113").Value = 349.5
$ws.Range("K113").Value = 349.5
$ws.Range("M113").Value = 1820.5
$ws.Range("H126").Value = 4501.2
$ws.Range("I126").Value = 4754.4443
$ws.Range("J126").Value = 4121.3335
$ws.Range("K126").Value = 14263.3329
$ws.Range("L126").Value = 12364.0005
$ws.Range("M126").Value = -11793.3329
$ws.Range("N126").Value = -17304.0005
$ws.Range("H134").Value = 2018.7333
$ws.Range("I134").Value = 1711.619
$ws.Range("K134").Value = 5134.857
$ws.Range("M134").Value = -2599.857

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 1203.8148
$ws.Range("I14").Value = 1203.8148
$ws.Range("K14").Value = 3611.4444
$ws.Range("M14").Value = -3438.4444
$ws.Range("H124").Value = 5429.4443
$ws.Range("I124").Value = 4414.5
$ws.Range("J124").Value = 5719.4287
$ws.Range("K124").Value = 13243.5
$ws.Range("L124").Value = 17158.2861
$ws.Range("M124").Value = -8333.5
$ws.Range("N124").Value = -26978.2861
$ws.Range("H129").Value = 3439.9285
$ws.Range("J129").Value = 3827.4167
$ws.Range("L129").Value = 11482.2501
$ws.Range("N129").Value = -21482.2501
$ws.Range("H137").Value = 4869.6284
$ws.Range("I137").Value = 664.2
$ws.Range("J137").Value = 5157.6714
$ws.Range("K137").Value = 1992.6
$ws.Range("L137").Value = 15473.0142
$ws.Range("M137").Value = 3107.4
$ws.Range("N137").Value = -25673.0142
$ws.Range("H138").Value = 5588.0713
$ws.Range("I138").Value = 3137.5557
$ws.Range("K138").Value = 9412.667099999999
$ws.Range("M138").Value = -4272.667099999999

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 183919.8
$ws.Range("I122").Value = 301666.34
$ws.Range("J122").Value = 7300
$ws.Range("K122").Value = 904999.02
$ws.Range("L122").Value = 21900
$ws.Range("M122").Value = -902549.02
$ws.Range("N122").Value = -26800
$ws.Range("H132").Value = 4200.1177
$ws.Range("I132").Value = 4036.5417
$ws.Range("K132").Value = 12109.6251
$ws.Range("M132").Value = -9579.625100000001

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 3432.524
$ws.Range("I122").Value = 3293.5405
$ws.Range("J122").Value = 4461
$ws.Range("K122").Value = 9880.621500000001
$ws.Range("L122").Value = 13383
$ws.Range("M122").Value = -7430.621500000001
$ws.Range("N122").Value = -18283
$ws.Range("H132").Value = 3911.6897
$ws.Range("I132").Value = 3231.52
$ws.Range("K132").Value = 9694.559999999999
$ws.Range("M132").Value = -7164.559999999999
$ws.Range("H136").Value = 4919.5884
$ws.Range("J136").Value = 5766.8184
$ws.Range("L136").Value = 17300.4552
$ws.Range("N136").Value = -22400.4552

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 889.10254
$ws.Range("I113").Value = 903
$ws.Range("K113").Value = 2709
$ws.Range("M113").Value = -539
$ws.Range("H122").Value = 2663.25
$ws.Range("I122").Value = 2213.6
$ws.Range("K122").Value = 6640.799999999999
$ws.Range("M122").Value = -4190.799999999999
$ws.Range("H132").Value = 282283.56
$ws.Range("I132").Value = 374308.28
$ws.Range("J132").Value = 6209.5
$ws.Range("K132").Value = 1122924.84
$ws.Range("L132").Value = 18628.5
$ws.Range("M132").Value = -1120394.84
$ws.Range("N132").Value = -23688.5
